# Auto-generated edit script: updates cached numeric values on the
# "Ragnarok_Profits" workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to match the latest scheduled-runner data pull.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 3: H3, J3, L3, N3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null
# Row 32: H32, I32, K32, M32
$ws.Range("H32").Value = 2989
$ws.Range("I32").Value = 2989
$ws.Range("K32").Value = 2989
$ws.Range("M32").Value = -2663
# Row 70: H70, J70, L70, N70
$ws.Range("H70").Value = 1527758.4
$ws.Range("J70").Value = 3849.4285
$ws.Range("L70").Value = 11548.2855
$ws.Range("N70").Value = -12088.2855
# Row 73: H73, J73, L73, N73
$ws.Range("H73").Value = 1527758.4
$ws.Range("J73").Value = 3849.4285
$ws.Range("L73").Value = 11548.2855
$ws.Range("N73").Value = -13420.2855
# Row 94: H94, I94, K94, M94
$ws.Range("H94").Value = 5332.25
$ws.Range("I94").Value = 5332.25
$ws.Range("K94").Value = 5332.25
$ws.Range("M94").Value = -4881.25
# Row 95: H95, J95, L95, N95
$ws.Range("H95").Value = 50640.25
$ws.Range("J95").Value = 50640.25
$ws.Range("L95").Value = 50640.25
$ws.Range("N95").Value = -56132.25
# Row 102: H102, J102, L102, N102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null
# Row 107: H107, I107, J107, K107, L107, M107, N107
$ws.Range("H107").Value = 1004.1875
$ws.Range("I107").Value = 580
$ws.Range("J107").Value = 1937.4
$ws.Range("K107").Value = 580
$ws.Range("L107").Value = 1937.4
$ws.Range("M107").Value = 1340
$ws.Range("N107").Value = -5777.4
# Row 111: H111, I111, J111, K111, L111, M111, N111
$ws.Range("H111").Value = 2750
$ws.Range("I111").Value = 2500
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 7500
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -4433
$ws.Range("N111").Value = -15134
# Row 113: H113, I113, K113, M113
$ws.Range("H113").Value = 8535.333000000001
$ws.Range("I113").Value = 7800
$ws.Range("K113").Value = 7800
$ws.Range("M113").Value = -4546
# Row 116: H116, I116, J116, K116, L116, M116, N116
$ws.Range("H116").Value = 8907.223
$ws.Range("I116").Value = 8377.4
$ws.Range("J116").Value = 9569.5
$ws.Range("K116").Value = 8377.4
$ws.Range("L116").Value = 9569.5
$ws.Range("M116").Value = -4935.4
$ws.Range("N116").Value = -16453.5
# Row 137: H137, I137, K137, M137
$ws.Range("H137").Value = 1857448.9
$ws.Range("I137").Value = 2443.6667
$ws.Range("K137").Value = 7331.000100000001
$ws.Range("M137").Value = -4781.000100000001

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 44: H44, I44, J44, K44, L44, M44, N44
$ws.Range("H44").Value = 52499.75
$ws.Range("I44").Value = 20000
$ws.Range("J44").Value = 63333
$ws.Range("K44").Value = 20000
$ws.Range("L44").Value = 63333
$ws.Range("M44").Value = -19512
$ws.Range("N44").Value = -64309
# Row 55: H55, I55, J55, K55, L55, M55, N55
$ws.Range("H55").Value = 33624.5
$ws.Range("I55").Value = 28166
$ws.Range("J55").Value = 50000
$ws.Range("K55").Value = 28166
$ws.Range("L55").Value = 50000
$ws.Range("M55").Value = -27851
$ws.Range("N55").Value = -50630
# Row 97: H97, I97, J97, K97, L97, M97, N97
$ws.Range("H97").Value = 965.8095
$ws.Range("I97").Value = 886.2820400000001
$ws.Range("J97").Value = 1999.6666
$ws.Range("K97").Value = 886.2820400000001
$ws.Range("L97").Value = 1999.6666
$ws.Range("M97").Value = -390.2820400000001
$ws.Range("N97").Value = -2991.6666
# Row 103: H103, J103, L103, N103
$ws.Range("H103").Value = 90000
$ws.Range("J103").Value = 90000
$ws.Range("L103").Value = 90000
$ws.Range("N103").Value = -92344
# Row 110: H110, I110, K110, M110
$ws.Range("H110").Value = 10385.286
$ws.Range("I110").Value = 10539.4
$ws.Range("K110").Value = 10539.4
$ws.Range("M110").Value = -8494.4
# Row 122: H122, I122, K122, M122
$ws.Range("H122").Value = 4113.2144
$ws.Range("I122").Value = 4072.889
$ws.Range("K122").Value = 12218.667
$ws.Range("M122").Value = -9768.667000000001
# Row 126: H126, I126, K126, M126
$ws.Range("H126").Value = 35714284
$ws.Range("I126").Value = 35714284
$ws.Range("K126").Value = 107142852
$ws.Range("M126").Value = -107140382

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 7: H7, I7, J7, K7, L7, M7, N7
$ws.Range("H7").Value = 3333696.2
$ws.Range("I7").Value = 544.5
$ws.Range("J7").Value = 10000000
$ws.Range("K7").Value = 544.5
$ws.Range("L7").Value = 10000000
$ws.Range("M7").Value = -431.5
$ws.Range("N7").Value = -10000226
# Row 94: H94, I94, J94, K94, L94, M94, N94
$ws.Range("H94").Value = 2221.7942
$ws.Range("I94").Value = 2109.24
$ws.Range("J94").Value = 2534.4443
$ws.Range("K94").Value = 2109.24
$ws.Range("L94").Value = 2534.4443
$ws.Range("M94").Value = -1658.24
$ws.Range("N94").Value = -3436.4443
# Row 107: H107, I107, J107, K107, L107, M107, N107
$ws.Range("H107").Value = 4221.1
$ws.Range("I107").Value = 4198.3125
$ws.Range("J107").Value = 4312.25
$ws.Range("K107").Value = 4198.3125
$ws.Range("L107").Value = 4312.25
$ws.Range("M107").Value = -2278.3125
$ws.Range("N107").Value = -8152.25

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 2: H2, I2, J2, K2, L2, M2, N2
$ws.Range("H2").Value = 3000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = -3226
# Row 16: H16, I16, J16, K16, L16, M16, N16
$ws.Range("H16").Value = 5269685
$ws.Range("I16").Value = 7697715.5
$ws.Range("J16").Value = 8952.166999999999
$ws.Range("K16").Value = 7697715.5
$ws.Range("L16").Value = 8952.166999999999
$ws.Range("M16").Value = -7697428.5
$ws.Range("N16").Value = -9526.166999999999
# Row 113: H113, I113, J113, K113, L113, M113, N113
$ws.Range("H113").Value = 5269685
$ws.Range("I113").Value = 7697715.5
$ws.Range("J113").Value = 8952.166999999999
$ws.Range("K113").Value = 7697715.5
$ws.Range("L113").Value = 8952.166999999999
$ws.Range("M113").Value = -7695545.5
$ws.Range("N113").Value = -13292.167
# Row 132: H132, I132, K132, M132
$ws.Range("H132").Value = 2755.2666
$ws.Range("I132").Value = 2729.0476
$ws.Range("K132").Value = 8187.1428
$ws.Range("M132").Value = -5657.1428

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 3: H3, I3, J3, K3, L3, M3, N3
$ws.Range("H3").Value = 1848.6
$ws.Range("I3").Value = 1848.6
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5545.799999999999
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -5433.799999999999
$ws.Range("N3").Value = $null
# Row 106: H106, J106, L106, N106
$ws.Range("H106").Value = 13131.3
$ws.Range("J106").Value = 20552.166
$ws.Range("L106").Value = 61656.49800000001
$ws.Range("N106").Value = -63548.49800000001
# Row 107: H107, J107, L107, N107
$ws.Range("H107").Value = 6072139.5
$ws.Range("J107").Value = 8278865
$ws.Range("L107").Value = 24836595
$ws.Range("N107").Value = -24840435
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 4013.9285
$ws.Range("I132").Value = 995
$ws.Range("J132").Value = 4246.154
$ws.Range("K132").Value = 8955
$ws.Range("L132").Value = 38215.38600000001
$ws.Range("M132").Value = -6425
$ws.Range("N132").Value = -43275.38600000001
# Row 136: H136, I136, K136, M136
$ws.Range("H136").Value = 7869.4443
$ws.Range("I136").Value = 1586.1666
$ws.Range("K136").Value = 4758.4998
$ws.Range("M136").Value = 341.5002000000004
# Row 141: H141, I141, K141, M141
$ws.Range("H141").Value = 7073.143
$ws.Range("I141").Value = 5053.154
$ws.Range("K141").Value = 15159.462
$ws.Range("M141").Value = -9979.462000000001

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 80: H80, I80, J80, K80, L80, M80, N80
$ws.Range("H80").Value = 2840.0557
$ws.Range("I80").Value = 1759.4445
$ws.Range("J80").Value = 3920.6667
$ws.Range("K80").Value = 1759.4445
$ws.Range("L80").Value = 3920.6667
$ws.Range("M80").Value = -761.4445000000001
$ws.Range("N80").Value = -5916.6667
# Row 83: H83, I83, J83, K83, L83, M83, N83
$ws.Range("H83").Value = 2840.0557
$ws.Range("I83").Value = 1759.4445
$ws.Range("J83").Value = 3920.6667
$ws.Range("K83").Value = 8797.2225
$ws.Range("L83").Value = 19603.3335
$ws.Range("M83").Value = -3805.2225
$ws.Range("N83").Value = -29587.3335
# Row 113: H113, J113, L113, N113
$ws.Range("H113").Value = 1237794.5
$ws.Range("J113").Value = 3706402.5
$ws.Range("L113").Value = 3706402.5
$ws.Range("N113").Value = -3710742.5
# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 3644.1177
$ws.Range("I122").Value = 3445.6667
$ws.Range("J122").Value = 4120.4
$ws.Range("K122").Value = 10337.0001
$ws.Range("L122").Value = 12361.2
$ws.Range("M122").Value = -7887.000100000001
$ws.Range("N122").Value = -17261.2

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16: H16, I16, J16, K16, L16, M16, N16
$ws.Range("H16").Value = 4204.4736
$ws.Range("I16").Value = 2688.8
$ws.Range("J16").Value = 5888.5557
$ws.Range("K16").Value = 2688.8
$ws.Range("L16").Value = 5888.5557
$ws.Range("M16").Value = -2518.8
$ws.Range("N16").Value = -6228.5557
# Row 106: H106, J106, L106, N106
$ws.Range("H106").Value = 20599.5
$ws.Range("J106").Value = 20599.5
$ws.Range("L106").Value = 20599.5
$ws.Range("N106").Value = -23123.5
# Row 122: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 3478.0425
$ws.Range("I122").Value = 3333.4092
$ws.Range("J122").Value = 5599.3335
$ws.Range("K122").Value = 10000.2276
$ws.Range("L122").Value = 16798.0005
$ws.Range("M122").Value = -7550.2276
$ws.Range("N122").Value = -21698.0005

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 15: H15, I15, J15, K15, L15, M15, N15
$ws.Range("H15").Value = 30282.143
$ws.Range("I15").Value = 2000
$ws.Range("J15").Value = 34995.832
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 34995.832
$ws.Range("M15").Value = -1712
$ws.Range("N15").Value = -35571.832
# Row 135: H135, J135, L135, N135
$ws.Range("H135").Value = 101175
$ws.Range("J135").Value = 101175
$ws.Range("L135").Value = 101175
$ws.Range("N135").Value = -111315
# Row 137: H137, J137, L137, N137
$ws.Range("H137").Value = 121119.2
$ws.Range("J137").Value = 135565.33
$ws.Range("L137").Value = 135565.33
$ws.Range("N137").Value = -145765.33

